# Clientnavigation.xlsx edit:
# - B1 value changes from "/clients" to "clients" (leading slash removed)
# - Active selection moves from C6 to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NavigateURL value for "createclient" (B1) by stripping the leading slash
$ws.Range("B1").Value = "clients"

# Move the selection to B2 to match the saved view state
$ws.Range("B2").Select()
